# Reorder/edit the BOM rows on the active sheet:
#  - Add C6 to the 10uF cap designator list and bump its quantity 2 -> 3
#  - Move the "Series 102 ... WR-TBL" connector row after the capacitor rows,
#    rename its designator from "Battery Connector, E-Match Connector" to
#    "E-Match, Power"
#  - Move the inductor (L1) row up, ahead of the transistor/LED rows
#  - Rename "LED?" -> "LED", "E-Match MOSFET" -> "MOSFET" (now on the
#    TRANS NPN row), and "Voltage Reg." -> "Reg."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: General Purpose Ceramic Capacitor, 1210, 10uF -> now C1, C5, C6 (qty 3)
$ws.Range("A4").Value = "General Purpose Ceramic Capacitor, 1210, 10uF, 10%, X7R, 15%, 25V"
$ws.Range("B4").Value = "C1, C5, C6"
$ws.Range("E4").Value = 3

# --- Row 5: General Purpose Ceramic Capacitor, 0805, 100nF -> C2 (qty 1)
$ws.Range("A5").Value = "General Purpose Ceramic Capacitor, 0805, 100nF, 5%, X7R, 0.15, 50V"
$ws.Range("B5").Value = "C2"
$ws.Range("E5").Value = 1

# --- Row 6: Cap Ceramic 22uF ... -> C3, C4 (qty 2)
$ws.Range("A6").Value = "Cap Ceramic 22uF 10V X7R ±20% SMD 1206 +125°C Embossed T/R"
$ws.Range("B6").Value = "C3, C4"
$ws.Range("E6").Value = 2

# --- Row 7: Series 102 connector -> E-Match, Power (qty 2, unchanged)
$ws.Range("A7").Value = "Series 102 - 5.00 mm Horizontal Entry Modular with Pressure Clamp WR-TBL, 2 pin"
$ws.Range("B7").Value = "E-Match, Power"

# --- Row 8: Shielded Power Inductor -> L1 (qty 1, unchanged)
$ws.Range("A8").Value = "Shielded Power Inductor WE-PD2SR, L=3.9 µH"
$ws.Range("B8").Value = "L1"

# --- Row 9: LED RED CLEAR CHIP SMD -> LED (qty 1, unchanged)
$ws.Range("A9").Value = "LED RED CLEAR CHIP SMD"
$ws.Range("B9").Value = "LED"

# --- Row 10: TRANS NPN 60V 1A SOT23-3 -> MOSFET (qty 1, unchanged)
$ws.Range("A10").Value = "TRANS NPN 60V 1A SOT23-3"
$ws.Range("B10").Value = "MOSFET"

# --- Row 13: designator Voltage Reg. -> Reg. (description/qty unchanged)
$ws.Range("B13").Value = "Reg."

# Re-stamp the quote-prefixed text style (preserved on every other Description/
# Designator cell) onto the cells we just rewrote, since setting .Value resets
# cell formatting to the default text style.
$ws.Range("A2").Copy()
$ws.Range("A4:A10").PasteSpecial(-4122)
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B4:B10").PasteSpecial(-4122)
$ws.Range("B13").PasteSpecial(-4122)

$excel.CutCopyMode = 0
